$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StepperPage")
$ws.Activate()

$ws.Range("D2").Value = "Automation"
$ws.Range("E2").Value = "Test"

$ws.Range("E2").Select()
